$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.057.24"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "'3.390.73"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'574.07"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").Value = "'137.59"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D8").Value = "'3.389.06"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").Value = "'0.388"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").Value = "'3.961.24"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").Value = "'26.21"
$ws.Range("E15").Value = "  +2.94%  "

$ws.Range("E16").Value = "  -2.42%  "

$ws.Range("D17").Value = "'3.386.46"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").Value = "'61.149.00"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").Value = "'14.05"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("D20").Value = "'5.83"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").Value = "'9.46"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").Value = "'377.34"
$ws.Range("E22").Value = "  -3.06%  "

$ws.Range("D23").Value = "'0.556"
$ws.Range("E23").Value = "  -2.65%  "

$ws.Range("D24").Value = "'3.525.23"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("E26").Value = "  -2.56%  "

$ws.Range("D27").Value = "'71.25"
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").Value = "'1.77"
$ws.Range("E28").Value = "  +11.95%  "

$ws.Range("E29").Value = "  +7.66%  "

$ws.Range("D30").Value = "'7.52"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "'8.18"
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'23.72"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").Value = "'5.24"
$ws.Range("E36").Value = "  -3.55%  "

$ws.Range("D37").Value = "'1.56"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  -1.71%  "

$ws.Range("D39").Value = "'163.99"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").Value = "'0.0764"
$ws.Range("E40").Value = "  -3.07%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "'0.776"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("E43").Value = "  -3.13%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.43"
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'41.67"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "'1.21"
$ws.Range("E46").Value = "  -2.02%  "

$ws.Range("D47").Value = "'24.22"
$ws.Range("E47").Value = "  -3.56%  "

$ws.Range("D48").Value = "'2.467.16"
$ws.Range("E48").Value = "  +4.21%  "

$ws.Range("D49").Value = "'23.25"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").Value = "'6.82"
$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("D51").Value = "'2.43"
$ws.Range("E51").Value = "  +4.45%  "
